# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.683.65'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '1.621.91'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.44'
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0611'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.35'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').Value = '1.852.11'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '1.628.98'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.97'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.512'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').Value = '26.711.40'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.95'
$ws.Range('E18').Value = '  +9.86%  '
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('D20').Value = '0.0₃0728'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.27'
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.10'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.83'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('E28').Value = '  +2.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.67'
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('E31').Value = '  +1.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('D33').Value = '1.464.93'
$ws.Range('E33').Value = '  +9.49%  '
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.569'
$ws.Range('E37').Value = '  -1.92%  '
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.839'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.95'
$ws.Range('E40').Value = '  +3.00%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.958'
$ws.Range('E42').Value = '  +3.02%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.21'
$ws.Range('E43').Value = '  +3.02%  '
$ws.Range('D44').Value = '1.763.62'
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.93'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.53'
$ws.Range('E47').Value = '  +3.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.50'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0504'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0965'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.50'
$ws.Range('E51').Value = '  +2.17%  '
